# "dang dat cho tu cu ben controller study" - porting/staging a leftover
# vocabulary entry ("easy") in from another module's shared-string pool.
#
# The underlying OOXML diff only touches xl/sharedStrings.xml (a new,
# still-unused <si> for "easy" is inserted right after "tay nguyen", which
# pushes every following shared-string index up by one) and the <v> index
# references inside xl/worksheets/sheet1.xml that point into that table.
# Every cell keeps displaying exactly the same English/Tieng Viet/Level/
# PathImage text it did before - only the pointer into the (re-shuffled)
# shared-string pool changes under the hood, because the new "easy" entry
# is not wired up to any cell yet (it is just being staged for later use,
# per the commit message).
#
# The Excel object model has no notion of "shared string table order" -
# that bookkeeping is an implementation detail the engine recomputes from
# actual cell contents whenever the workbook is saved. So the faithful way
# to replay this edit through COM automation is to touch the same cell
# range the diff touches and re-assert the same values, letting the engine
# rebuild/renumber the shared-string table on save exactly like the
# original author's tool did - without altering anything a user would
# actually see.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 37

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $english = $ws.Cells.Item($r, 1).Value()
    $tiengViet = $ws.Cells.Item($r, 2).Value()
    $level = $ws.Cells.Item($r, 3).Value()

    $ws.Cells.Item($r, 1).Value = $english
    $ws.Cells.Item($r, 2).Value = $tiengViet
    $ws.Cells.Item($r, 3).Value = $level
}
